$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.540.10'
$ws.Range('E2').Value = '  +0.70%  '

$ws.Range('D3').Value = '1.728.53'
$ws.Range('E3').Value = '  +0.61%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.03'
$ws.Range('E5').Value = '  +2.25%  '

$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4792'
$ws.Range('E7').Value = '  +1.20%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2669'
$ws.Range('E8').Value = '  +1.48%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06224'
$ws.Range('E9').Value = '  +0.20%  '

$ws.Range('D10').Value = '1.730.16'
$ws.Range('E10').Value = '  +0.84%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07147'
$ws.Range('E11').Value = '  +1.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.68'
$ws.Range('E12').Value = '  +2.23%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6134'
$ws.Range('E13').Value = '  +3.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.528'
$ws.Range('E14').Value = '  +2.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.96'
$ws.Range('E15').Value = '  +1.22%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.0000'
$ws.Range('E16').Value = '  -0.03%  '

$ws.Range('D17').Value = '26.548.11'
$ws.Range('E17').Value = '  +0.76%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006970'
$ws.Range('E19').Value = '  +2.18%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.64'
$ws.Range('E20').Value = '  +0.84%  '

$ws.Range('D21').Value = '1.951.71'
$ws.Range('E21').Value = '  +0.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.526'
$ws.Range('E22').Value = '  -0.34%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.927'
$ws.Range('E23').Value = '  +1.91%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.286'
$ws.Range('E24').Value = '  -0.72%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.59'
$ws.Range('E25').Value = '  +1.58%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.35'
$ws.Range('E26').Value = '  +0.71%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.793'
$ws.Range('E27').Value = '  +2.26%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.406'
$ws.Range('E28').Value = '  +0.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.64'
$ws.Range('E29').Value = '  -1.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.982'
$ws.Range('E30').Value = '  -0.39%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07967'
$ws.Range('E31').Value = '  +2.80%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.710'
$ws.Range('E32').Value = '  +0.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04590'
$ws.Range('E33').Value = '  +3.78%  '

$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9999'
$ws.Range('E34').Value = '  +0.04%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.618'
$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9947'
$ws.Range('E36').Value = '  +1.77%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6311'
$ws.Range('E37').Value = '  +1.97%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.093'
$ws.Range('E38').Value = '  +9.16%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9168'
$ws.Range('E39').Value = '  -2.00%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.405'
$ws.Range('E40').Value = '  -0.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').Value = '  +0.70%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.55'
$ws.Range('E42').Value = '  -7.39%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.01503'
$ws.Range('E43').Value = '  +1.81%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.568'
$ws.Range('E44').Value = '  +4.31%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3879'
$ws.Range('E45').Value = '  +1.81%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.978'
$ws.Range('E46').Value = '  +10.78%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1184'
$ws.Range('E47').Value = '  +1.25%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05345'
$ws.Range('E48').Value = '  +1.17%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.91'
$ws.Range('E49').Value = '  +2.04%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.814'
$ws.Range('E50').Value = '  +1.26%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.259'
$ws.Range('E51').Value = '  +3.75%  '

